# Apply the per-row cryptocurrency price (D) and 1h volume % change (E)
# updates described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "23.183.33"
$ws.Range("E2").Value = "  +0.31%  "
$ws.Range("D3").Value = "1.601.74"
$ws.Range("E3").Value = "  -0.01%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D6").Value = "'303.24"
$ws.Range("E6").Value = "  +0.53%  "
$ws.Range("D7").Value = "'0.3779"
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("D8").Value = "'52.06"
$ws.Range("E8").Value = "  +4.36%  "
$ws.Range("D9").Value = "'0.3615"
$ws.Range("E10").Value = "  +0.20%  "
$ws.Range("E11").Value = "  +0.08%  "
$ws.Range("D12").Value = "'0.08117"
$ws.Range("E12").Value = "  -0.62%  "
$ws.Range("E13").Value = "  -1.17%  "
$ws.Range("D14").Value = "'6.586"
$ws.Range("E14").Value = "  -0.13%  "
$ws.Range("E15").Value = "  +0.25%  "
$ws.Range("D16").Value = "'0.00001244"
$ws.Range("E16").Value = "  -1.33%  "
$ws.Range("D17").Value = "1.602.08"
$ws.Range("E17").Value = "  +0.20%  "
$ws.Range("D18").Value = "'94.10"
$ws.Range("E18").Value = "  +2.46%  "
$ws.Range("D19").Value = "'0.06874"
$ws.Range("E19").Value = "  +0.30%  "
$ws.Range("D20").Value = "'18.07"
$ws.Range("E20").Value = "  -0.89%  "
$ws.Range("D21").Value = "'6.547"
$ws.Range("E21").Value = "  -0.70%  "
$ws.Range("E22").Value = "  -0.06%  "
$ws.Range("D23").Value = "'12.98"
$ws.Range("E23").Value = "  -0.52%  "
$ws.Range("D24").Value = "23.190.26"
$ws.Range("E24").Value = "  +0.37%  "
$ws.Range("D25").Value = "'2.401"
$ws.Range("E25").Value = "  +2.49%  "
$ws.Range("D26").Value = "'2.977"
$ws.Range("E26").Value = "  +8.68%  "
$ws.Range("E27").Value = "  +0.46%  "
$ws.Range("D28").Value = "'149.37"
$ws.Range("E28").Value = "  -0.31%  "
$ws.Range("D29").Value = "'5.255"
$ws.Range("E29").Value = "  -0.16%  "
$ws.Range("D30").Value = "'133.96"
$ws.Range("E30").Value = "  +1.04%  "
$ws.Range("D31").Value = "'2.379"
$ws.Range("E31").Value = "  -0.56%  "
$ws.Range("D32").Value = "'6.776"
$ws.Range("E32").Value = "  -1.01%  "
$ws.Range("D33").Value = "1.778.49"
$ws.Range("E33").Value = "  +0.06%  "
$ws.Range("D34").Value = "'0.9695"
$ws.Range("E34").Value = "  +1.48%  "
$ws.Range("D35").Value = "'0.07517"
$ws.Range("E35").Value = "  -2.32%  "
$ws.Range("D36").Value = "'10.31"
$ws.Range("E36").Value = "  +2.31%  "
$ws.Range("D37").Value = "'0.02720"
$ws.Range("E37").Value = "  -0.03%  "
$ws.Range("E38").Value = "  -2.04%  "
$ws.Range("D39").Value = "'0.08806"
$ws.Range("E39").Value = "  -1.02%  "
$ws.Range("D40").Value = "'6.076"
$ws.Range("E40").Value = "  -3.20%  "
$ws.Range("D41").Value = "'0.7111"
$ws.Range("E41").Value = "  +0.39%  "
$ws.Range("D42").Value = "'1.359"
$ws.Range("E42").Value = "  -0.85%  "
$ws.Range("D43").Value = "'12.50"
$ws.Range("E43").Value = "  -0.88%  "
$ws.Range("D44").Value = "'15.67"
$ws.Range("E44").Value = "  +1.91%  "
$ws.Range("D45").Value = "'0.6527"
$ws.Range("D46").Value = "'2.312"
$ws.Range("E46").Value = "  -0.29%  "
$ws.Range("D47").Value = "'4.018"
$ws.Range("E47").Value = "  +0.53%  "
$ws.Range("D48").Value = "'132.16"
$ws.Range("E48").Value = "  +0.60%  "
$ws.Range("D49").Value = "'0.07957"
$ws.Range("E49").Value = "  +0.21%  "
$ws.Range("D50").Value = "'1.200"
$ws.Range("E50").Value = "  -2.81%  "
$ws.Range("D51").Value = "'1.214"
$ws.Range("E51").Value = "  +1.47%  "
